# #5: property boat&car done
# Adds the missing property/category/date/legislator columns (H:N) to the
# "汽車" (car) sheet, matching the layout already used by the other
# property sheets, and fixes the row-2/row-3 name & register_date values
# that had been mixed up between the two cars.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Header row (row 1) --------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# New header cells (H1:N1) need the same bold / centred / thin-border look
# already used by the other header cells on this row (style index 1).
$headerRange = $ws.Range("H1:N1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Row 2 (index 35, TOYOTAPREVIA) --------------------------------------
$ws.Range("B2").Value = "TOYOTAPREVIA"
$ws.Range("E2").Value = "100年10月25曰"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-11-21"
$ws.Range("K2").Value = "孫大千"
$ws.Range("L2").Value = 919
$ws.Range("M2").Value = "tmpc6841"
$ws.Range("N2").Value = 35

# ---- Row 3 (index 36, 納智捷G91SPCA) -------------------------------------
$ws.Range("B3").Value = "納智捷G91SPCA"
$ws.Range("E3").Value = "100年04月28日"
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2011-11-21"
$ws.Range("K3").Value = "孫大千"
$ws.Range("L3").Value = 919
$ws.Range("M3").Value = "tmpc6841"
$ws.Range("N3").Value = 36
